# Update "Daily Data 2024-12-05 TODAY.xlsx" with refreshed source-feed
# values (DKIS INFO, ASX, 7 Day Outlook) and bump the DKIS date header.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# DKIS INFO sheet
# ---------------------------------------------------------------
$wsDKIS = $wb.Worksheets.Item("DKIS INFO")

$wsDKIS.Range("A1").Value = "04 Dec 2024"

$wsDKIS.Range("A2").Value = 135.21
$wsDKIS.Range("B2").Value = 130.26
$wsDKIS.Range("C2").Value = 127.77
$wsDKIS.Range("D2").Value = 110.63
$wsDKIS.Range("E2").Value = 132.42
$wsDKIS.Range("F2").Value = 132.28
$wsDKIS.Range("G2").Value = 168.45
$wsDKIS.Range("H2").Value = 171
$wsDKIS.Range("I2").Value = 80.66
$wsDKIS.Range("J2").Value = 82.48999999999999
$wsDKIS.Range("K2").Value = 151.78
$wsDKIS.Range("L2").Value = 135.21
$wsDKIS.Range("M2").Value = 246.1400146484375

# ---------------------------------------------------------------
# ASX sheet
# ---------------------------------------------------------------
$wsASX = $wb.Worksheets.Item("ASX")

$wsASX.Range("B2").Value = 129.54
$wsASX.Range("C2").Value = 121.16

$wsASX.Range("B3").Value = 113.12
$wsASX.Range("C3").Value = 102.06

$wsASX.Range("B4").Value = 108.3
$wsASX.Range("C4").Value = 102.66

$wsASX.Range("B5").Value = 80.33
$wsASX.Range("C5").Value = 73.52

# ---------------------------------------------------------------
# 7 Day Outlook sheet (values stored as text, like the source data)
# A leading apostrophe tells Excel to keep the numeric-looking
# value as text instead of auto-converting it to a number.
# ---------------------------------------------------------------
$wsOutlook = $wb.Worksheets.Item("7 Day Outlook")

# NSW1
$wsOutlook.Range("B2").Value = "'10833"
$wsOutlook.Range("C2").Value = "'11158"
$wsOutlook.Range("D2").Value = "'8865"
$wsOutlook.Range("E2").Value = "'9052"
$wsOutlook.Range("F2").Value = "'9415"
$wsOutlook.Range("G2").Value = "'9112"
$wsOutlook.Range("H2").Value = "'8900"
$wsOutlook.Range("I2").Value = "'1950"
$wsOutlook.Range("J2").Value = "'1923"
$wsOutlook.Range("K2").Value = "'3608"
$wsOutlook.Range("L2").Value = "'3268"
$wsOutlook.Range("M2").Value = "'2830"
$wsOutlook.Range("N2").Value = "'3687"
$wsOutlook.Range("O2").Value = "'3935"

# QLD1
$wsOutlook.Range("B3").Value = "'8380"
$wsOutlook.Range("C3").Value = "'8180"
$wsOutlook.Range("D3").Value = "'9096"
$wsOutlook.Range("E3").Value = "'8987"
$wsOutlook.Range("F3").Value = "'8841"
$wsOutlook.Range("G3").Value = "'8893"
$wsOutlook.Range("H3").Value = "'8835"
$wsOutlook.Range("I3").Value = "'2812"
$wsOutlook.Range("J3").Value = "'3238"
$wsOutlook.Range("K3").Value = "'2629"
$wsOutlook.Range("L3").Value = "'2911"
$wsOutlook.Range("M3").Value = "'2865"
$wsOutlook.Range("N3").Value = "'2822"
$wsOutlook.Range("O3").Value = "'3113"

# SA1
$wsOutlook.Range("B4").Value = "'2361"
$wsOutlook.Range("C4").Value = "'1605"
$wsOutlook.Range("D4").Value = "'1403"
$wsOutlook.Range("E4").Value = "'1592"
$wsOutlook.Range("F4").Value = "'1585"
$wsOutlook.Range("G4").Value = "'1453"
$wsOutlook.Range("H4").Value = "'1501"
$wsOutlook.Range("I4").Value = "'979"
$wsOutlook.Range("J4").Value = "'1328"
$wsOutlook.Range("K4").Value = "'1982"
$wsOutlook.Range("L4").Value = "'1884"
$wsOutlook.Range("M4").Value = "'2247"
$wsOutlook.Range("N4").Value = "'2208"
$wsOutlook.Range("O4").Value = "'2067"

# VIC1
$wsOutlook.Range("B5").Value = "'7011"
$wsOutlook.Range("C5").Value = "'5694"
$wsOutlook.Range("D5").Value = "'4762"
$wsOutlook.Range("E5").Value = "'5210"
$wsOutlook.Range("F5").Value = "'5238"
$wsOutlook.Range("G5").Value = "'5108"
$wsOutlook.Range("H5").Value = "'5183"
$wsOutlook.Range("I5").Value = "'2607"
$wsOutlook.Range("J5").Value = "'4283"
$wsOutlook.Range("K5").Value = "'4591"
$wsOutlook.Range("L5").Value = "'3975"
$wsOutlook.Range("M5").Value = "'4051"
$wsOutlook.Range("N5").Value = "'4905"
$wsOutlook.Range("O5").Value = "'4656"

# TAS1
$wsOutlook.Range("B6").Value = "'1209"
$wsOutlook.Range("C6").Value = "'1195"
$wsOutlook.Range("D6").Value = "'1125"
$wsOutlook.Range("E6").Value = "'1159"
$wsOutlook.Range("F6").Value = "'1183"
$wsOutlook.Range("G6").Value = "'1177"
$wsOutlook.Range("H6").Value = "'1172"
$wsOutlook.Range("I6").Value = "'939"
$wsOutlook.Range("J6").Value = "'1042"
$wsOutlook.Range("K6").Value = "'1405"
$wsOutlook.Range("L6").Value = "'1292"
$wsOutlook.Range("M6").Value = "'1193"
$wsOutlook.Range("N6").Value = "'1049"
$wsOutlook.Range("O6").Value = "'1150"
